$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "288.28"
Set-TextValue $ws "E2" "0.97%"
Set-TextValue $ws "D3" "29.30"
Set-TextValue $ws "E3" "2.48%"
Set-TextValue $ws "D4" "5.095"
Set-TextValue $ws "E4" "2.42%"
Set-TextValue $ws "D5" "0.06683"
Set-TextValue $ws "E5" "2.99%"
Set-TextValue $ws "D6" "7.322"
Set-TextValue $ws "E6" "1.32%"
Set-TextValue $ws "D7" "3.407"
Set-TextValue $ws "E7" "1.18%"
Set-TextValue $ws "D8" "1.369"
Set-TextValue $ws "E8" "1.79%"
Set-TextValue $ws "D9" "0.9179"
Set-TextValue $ws "E9" "0.55%"
Set-TextValue $ws "D10" "0.1589"
Set-TextValue $ws "E10" "3.28%"
Set-TextValue $ws "D11" "0.06758"
Set-TextValue $ws "E11" "7.93%"
Set-TextValue $ws "D12" "0.07585"
Set-TextValue $ws "E12" "-0.76%"
Set-TextValue $ws "D13" "0.02936"
Set-TextValue $ws "E13" "-1.68%"
Set-TextValue $ws "D14" "0.08979"
Set-TextValue $ws "E14" "0.27%"
Set-TextValue $ws "D15" "0.001584"
Set-TextValue $ws "E15" "-0.77%"
Set-TextValue $ws "D16" "0.04510"
Set-TextValue $ws "E16" "1.24%"
Set-TextValue $ws "D17" "0.0006480"
Set-TextValue $ws "E17" "-1.14%"
Set-TextValue $ws "D18" "0.006286"
Set-TextValue $ws "E18" "4.33%"
Set-TextValue $ws "D19" "3.446"
Set-TextValue $ws "E19" "-0.45%"
Set-TextValue $ws "D20" "2.222"
Set-TextValue $ws "E20" "-0.93%"
Set-TextValue $ws "E21" "2.00%"
Set-TextValue $ws "E22" "-2.47%"
Set-TextValue $ws "D23" "4.076"
Set-TextValue $ws "E23" "2.38%"
Set-TextValue $ws "E24" "1.80%"
Set-TextValue $ws "E25" "0.22%"
Set-TextValue $ws "D26" "0.004112"
Set-TextValue $ws "E26" "-4.92%"
Set-TextValue $ws "E27" "1.71%"
Set-TextValue $ws "D28" "0.0001617"
Set-TextValue $ws "E28" "-1.16%"
Set-TextValue $ws "D40" "0.04246"
Set-TextValue $ws "E40" "2.40%"
Set-TextValue $ws "D41" "0.006716"
Set-TextValue $ws "E41" "-0.17%"
Set-TextValue $ws "D42" "0.1238"
Set-TextValue $ws "E42" "0.53%"
Set-TextValue $ws "E43" "5.15%"
Set-TextValue $ws "D44" "0.01343"
Set-TextValue $ws "E44" "14.18%"
Set-TextValue $ws "D45" "0.00005711"
Set-TextValue $ws "E45" "5.95%"
Set-TextValue $ws "E46" "-3.59%"
Set-TextValue $ws "D47" "0.01306"
Set-TextValue $ws "E47" "-29.42%"
